$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E) for the account-statement data rows (16-22) is
# reordered: the most recent period (2306) moves to the top of the list and
# the oldest (2212) moves to the bottom (rows reversed).
#
#   Before: 2212, 2301, 2302, 2303, 2304, 2305, 2306
#   After : 2306, 2305, 2304, 2303, 2302, 2301, 2212
#
# The "Valor Mora" value of 35200 (previously tied to period 2306 / row 22)
# now travels with that period to row 16, while row 22 (now period 2212)
# gets the 60000 value that the other periods use.

$ws.Range("E16").Value = "2306"
$ws.Range("E17").Value = "2305"
$ws.Range("E18").Value = "2304"
$ws.Range("E19").Value = "2303"
$ws.Range("E20").Value = "2302"
$ws.Range("E21").Value = "2301"
$ws.Range("E22").Value = "2212"

$ws.Range("F16").Value = 35200
$ws.Range("F22").Value = 60000
